$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing row updated: sending cluster ECs -> Muc2/Agr2 -> ECs)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Muc2"
$ws.Range("C2").Value = "Agr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05442466666666667
$ws.Range("H2").Value = 0.163274
$ws.Range("I2").Value = 0.02819211748574673
$ws.Range("J2").Value = 0.02819211748574673
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.362105333333334
$ws.Range("N2").Value = 10.086316
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1829814620648889
$ws.Range("R2").Value = 1.646833158584
$ws.Range("S2").Value = 0.02819211748574673
$ws.Range("T2").Value = 0.02819211748574673

# Row 3 (existing row updated: sending cluster FAPs -> Muc2/Agr2 -> ECs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Muc2"
$ws.Range("C3").Value = "Agr2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.270304
$ws.Range("H3").Value = 3.810912
$ws.Range("I3").Value = 0.658020743240455
$ws.Range("J3").Value = 0.658020743240455
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.362105333333334
$ws.Range("N3").Value = 10.086316
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 4.270895853354667
$ws.Range("R3").Value = 38.438062680192
$ws.Range("S3").Value = 0.658020743240455
$ws.Range("T3").Value = 0.658020743240455

# Row 4 (new row: sending cluster M2 -> Muc2/Agr2 -> ECs)
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Muc2"
$ws.Range("C4").Value = "Agr2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5813196666666667
$ws.Range("H4").Value = 1.743959
$ws.Range("I4").Value = 0.3011250843265026
$ws.Range("J4").Value = 0.3011250843265026
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.362105333333334
$ws.Range("N4").Value = 10.086316
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.954457951671556
$ws.Range("R4").Value = 17.590121565044
$ws.Range("S4").Value = 0.3011250843265026
$ws.Range("T4").Value = 0.3011250843265026

# Row 5 (new row: sending cluster sCs -> Muc2/Agr2 -> ECs)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Muc2"
$ws.Range("C5").Value = "Agr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.024444
$ws.Range("H5").Value = 0.07333199999999999
$ws.Range("I5").Value = 0.01266205494729583
$ws.Range("J5").Value = 0.01266205494729583
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.362105333333334
$ws.Range("N5").Value = 10.086316
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.082183302768
$ws.Range("R5").Value = 0.739649724912
$ws.Range("S5").Value = 0.01266205494729583
$ws.Range("T5").Value = 0.01266205494729583
